# Update the "ATDD Scenarios" sheet sample data:
#  - Feature column (A2:A14): "Feature" -> "MyFeature"
#  - Sub Feature column (B2:B14): "SubFeature" -> "MySubFeature"
#  - Given-When-Then (Description) column (G) placeholder letters -> real sample text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

$ws.Range("A2:A14").Value = "MyFeature"
$ws.Range("B2:B14").Value = "MySubFeature"

$ws.Range("G4").Value  = "Some record"
$ws.Range("G5").Value  = "Do something"
$ws.Range("G6").Value  = "Something happens"
$ws.Range("G8").Value  = "Another record"
$ws.Range("G9").Value  = "Do something else"
$ws.Range("G10").Value = "An error was thrown"
$ws.Range("G12").Value = "Some other Record"
$ws.Range("G13").Value = "Do something in a page"
$ws.Range("G14").Value = "Something else happens"
